$d = $word.ActiveDocument

$find = "havainnointijaksot vuonna Kaksosten tähtikuvio 2022: 14.-23.2., 14.-24.3"
$replace = "Kaksosten tähtikuvio havainnointijaksot vuonna 2022: 14.-23.2., 14.-24.3"

$range = $d.Content
$range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
